$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 216.1448186666667
$ws.Range("H2").Value = 648.434456
$ws.Range("I2").Value = 0.9739197284603751
$ws.Range("J2").Value = 0.9739197284603752
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 24.030349
$ws.Range("N2").Value = 72.091047
$ws.Range("O2").Value = 0.8655386453070183
$ws.Range("P2").Value = 0.8655386453070184
$ws.Range("Q2").Value = 5194.035427101715
$ws.Range("R2").Value = 46746.31884391543
$ws.Range("S2").Value = 0.8429651624093721
$ws.Range("T2").Value = 0.8429651624093724

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 216.1448186666667
$ws.Range("H3").Value = 648.434456
$ws.Range("I3").Value = 0.9739197284603751
$ws.Range("J3").Value = 0.9739197284603752
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.12343
$ws.Range("N3").Value = 9.370290000000001
$ws.Range("O3").Value = 0.1125014610029717
$ws.Range("P3").Value = 0.1125014610029717
$ws.Range("Q3").Value = 675.1132109680267
$ws.Range("R3").Value = 6076.01889871224
$ws.Range("S3").Value = 0.1095673923514097
$ws.Range("T3").Value = 0.1095673923514097

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 216.1448186666667
$ws.Range("H4").Value = 648.434456
$ws.Range("I4").Value = 0.9739197284603751
$ws.Range("J4").Value = 0.9739197284603752
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.6096826666666667
$ws.Range("N4").Value = 1.829048
$ws.Range("O4").Value = 0.02195989369000996
$ws.Range("P4").Value = 0.02195989369000996
$ws.Range("Q4").Value = 131.7797494308764
$ws.Range("R4").Value = 1186.017744877888
$ws.Range("S4").Value = 0.0213871736995932
$ws.Range("T4").Value = 0.02138717369959321

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.181559666666666
$ws.Range("H5").Value = 9.544678999999999
$ws.Range("I5").Value = 0.01433568357434949
$ws.Range("J5").Value = 0.01433568357434949
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 24.030349
$ws.Range("N5").Value = 72.091047
$ws.Range("O5").Value = 0.8655386453070183
$ws.Range("P5").Value = 0.8655386453070184
$ws.Range("Q5").Value = 76.45398915432366
$ws.Range("R5").Value = 688.0859023889129
$ws.Range("S5").Value = 0.01240808814049253
$ws.Range("T5").Value = 0.01240808814049253

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.181559666666666
$ws.Range("H6").Value = 9.544678999999999
$ws.Range("I6").Value = 0.01433568357434949
$ws.Range("J6").Value = 0.01433568357434949
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.12343
$ws.Range("N6").Value = 9.370290000000001
$ws.Range("O6").Value = 0.1125014610029717
$ws.Range("P6").Value = 0.1125014610029717
$ws.Range("Q6").Value = 9.937378909656667
$ws.Range("R6").Value = 89.43641018691
$ws.Range("S6").Value = 0.001612785346590621
$ws.Range("T6").Value = 0.001612785346590621

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.181559666666666
$ws.Range("H7").Value = 9.544678999999999
$ws.Range("I7").Value = 0.01433568357434949
$ws.Range("J7").Value = 0.01433568357434949
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.6096826666666667
$ws.Range("N7").Value = 1.829048
$ws.Range("O7").Value = 0.02195989369000996
$ws.Range("P7").Value = 0.02195989369000996
$ws.Range("Q7").Value = 1.939741781732444
$ws.Range("R7").Value = 17.457676035592
$ws.Range("S7").Value = 0.0003148100872663367
$ws.Range("T7").Value = 0.0003148100872663367

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.606510333333333
$ws.Range("H8").Value = 7.819531
$ws.Range("I8").Value = 0.01174458796527538
$ws.Range("J8").Value = 0.01174458796527538
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 24.030349
$ws.Range("N8").Value = 72.091047
$ws.Range("O8").Value = 0.8655386453070183
$ws.Range("P8").Value = 0.8655386453070184
$ws.Range("Q8").Value = 62.63535298210633
$ws.Range("R8").Value = 563.718176838957
$ws.Range("S8").Value = 0.01016539475715356
$ws.Range("T8").Value = 0.01016539475715356

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.606510333333333
$ws.Range("H9").Value = 7.819531
$ws.Range("I9").Value = 0.01174458796527538
$ws.Range("J9").Value = 0.01174458796527538
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.12343
$ws.Range("N9").Value = 9.370290000000001
$ws.Range("O9").Value = 0.1125014610029717
$ws.Range("P9").Value = 0.1125014610029717
$ws.Range("Q9").Value = 8.141252570443333
$ws.Range("R9").Value = 73.27127313399001
$ws.Range("S9").Value = 0.001321283304971399
$ws.Range("T9").Value = 0.001321283304971399

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.606510333333333
$ws.Range("H10").Value = 7.819531
$ws.Range("I10").Value = 0.01174458796527538
$ws.Range("J10").Value = 0.01174458796527538
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.6096826666666667
$ws.Range("N10").Value = 1.829048
$ws.Range("O10").Value = 0.02195989369000996
$ws.Range("P10").Value = 0.02195989369000996
$ws.Range("Q10").Value = 1.589144170720889
$ws.Range("R10").Value = 14.302297536488
$ws.Range("S10").Value = 0.0002579099031504176
$ws.Range("T10").Value = 0.0002579099031504176

